$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (per-correct-answer mark) and the "Total" row
# (total correct marks + the "correct/max" summary string).
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
